# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row=3;  Value=319},
    @{Row=4;  Value=414},
    @{Row=5;  Value=1546},
    @{Row=7;  Value=2147},
    @{Row=8;  Value=3},
    @{Row=9;  Value=281},
    @{Row=11; Value=4797},
    @{Row=12; Value=9},
    @{Row=14; Value=299},
    @{Row=16; Value=25},
    @{Row=17; Value=167},
    @{Row=21; Value=3727},
    @{Row=22; Value=689},
    @{Row=23; Value=614},
    @{Row=27; Value=113},
    @{Row=29; Value=12},
    @{Row=34; Value=852},
    @{Row=35; Value=2336},
    @{Row=36; Value=421}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$updates4 = @(
    @{Row=3;  Value=319},
    @{Row=4;  Value=414},
    @{Row=5;  Value=1546},
    @{Row=7;  Value=2147},
    @{Row=8;  Value=3},
    @{Row=9;  Value=281},
    @{Row=11; Value=4797},
    @{Row=12; Value=9},
    @{Row=14; Value=299},
    @{Row=16; Value=25},
    @{Row=17; Value=167},
    @{Row=21; Value=3727},
    @{Row=22; Value=689},
    @{Row=23; Value=614},
    @{Row=27; Value=113},
    @{Row=29; Value=12},
    @{Row=35; Value=852},
    @{Row=36; Value=2336},
    @{Row=37; Value=421}
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates4) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
